$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# Change 1: "September 3, 2023"  ->  "September 15, 2023"   (submission deadline)
# Final text is split into three runs with identical (Calibri/20) formatting:
#   "September " | "15" | ", 2023"
# -------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("September 3, 2023")
if ($r1.Find.Found) {
    $s1 = $r1.Start

    # Update the day number in place first ("3" -> "15"); do this before any
    # run-isolation below so the whole run is still homogeneous when its text
    # changes (text replacement on a run boundary tends to re-merge sibling
    # runs sharing identical formatting).
    $d.Range($s1 + 10, $s1 + 11).Text = "15"

    # Final string is now "September 15, 2023" (19 chars):
    #   [0,10)  "September "
    #   [10,12) "15"
    #   [12,19) ", 2023"
    # Nudge formatting off/on (no visible effect) on each segment, right to
    # left, to force Word to keep them as separate runs instead of folding
    # everything back into a single run.
    $seg = $d.Range($s1 + 12, $s1 + 19)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0

    $seg = $d.Range($s1 + 10, $s1 + 12)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0

    $seg = $d.Range($s1, $s1 + 10)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0
}

# -------------------------------------------------------------------------
# Change 2: "September 25, 2023"  ->  "October 6, 2023"   (notification deadline)
# Final text is split into four runs with identical (Calibri/20) formatting:
#   "October" | " " | "6" | ", 2023"
# -------------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("September 25, 2023")
if ($r2.Find.Found) {
    $s2 = $r2.Start

    # Apply all text edits first (right to left so earlier offsets stay valid).
    $d.Range($s2 + 10, $s2 + 12).Text = "6"
    $d.Range($s2, $s2 + 9).Text = "October"

    # Final string is now "October 6, 2023" (15 chars):
    #   [0,7)  "October"
    #   [7,8)  " "
    #   [8,9)  "6"
    #   [9,15) ", 2023"
    $seg = $d.Range($s2 + 9, $s2 + 15)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0

    $seg = $d.Range($s2 + 8, $s2 + 9)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0

    $seg = $d.Range($s2 + 7, $s2 + 8)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0

    $seg = $d.Range($s2, $s2 + 7)
    $seg.Font.Bold = 1
    $seg.Font.Bold = 0
}
